# Add newly-logged Strava workouts (rows 87-92) to the weekly scoreboard,
# then extend the table's AutoFilter / _FilterDatabase range to match and
# move the selection to the new first empty row, mirroring what Excel does
# after a paste-append of new rows at the bottom of a filtered range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data rows -----------------------------------------------------
# Columns: A Participant | B Date | C Workout | D Total Duration |
#          E Total Distance | F Total Elevation | G Zone1 | H Zone2 |
#          I Zone3 | J Zone4 | K Zone5 | L Workout Type | M Week

$newRows = @(
    @{ Row=87; A="Jeremiah"; B=45464; C="Walk";    D=57;  E=2.66; F=171; G=45; H=1;  I=0;  J=0; K=0; L="Agile Antelope"; M=2 },
    @{ Row=88; A="Jeremiah"; B=45465; C="Run";     D=23;  E=2.29; F=207; G=0;  H=15; I=3;  J=0; K=0; L="Agile Antelope"; M=2 },
    @{ Row=89; A="Steven";   B=45465; C="Workout"; D=37;  E=0;    F=0;   G=8;  H=11; I=16; J=1; K=0; L="Agile Antelope"; M=2 },
    @{ Row=90; A="Steven";   B=45465; C="Run";     D=12;  E=1.04; F=46;  G=0;  H=2;  I=10; J=0; K=0; L="Agile Antelope"; M=2 },
    @{ Row=91; A="Eric";     B=45465; C="Workout"; D=69;  E=0;    F=0;   G=29; H=36; I=4;  J=0; K=0; L="Brave Leopard";  M=2 },
    @{ Row=92; A="Steven";   B=45465; C="Walk";    D=148; E=6.65; F=768; G=140;H=7;  I=2;  J=0; K=0; L="Agile Antelope"; M=2 }
)

# Copy the date format (style) already used by column B down onto the new
# date cells before writing values, so B87:B92 pick up the existing
# m/d/yyyy cell style (same xf as the rest of the column) instead of a
# brand-new number format.
$ws.Range("B2").Copy()
$ws.Range("B87:B92").PasteSpecial(-4122)
$excel.CutCopyMode = $false

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
}

# --- Extend the AutoFilter range to cover the new rows ------------------
# Range.AutoFilter() toggles the filter off if one is already active on
# that range, so call it twice: once to clear the existing A1:M71 filter,
# once more to (re)apply it across the full, now-larger A1:M92 table.
$ws.Range("A1:M92").AutoFilter() | Out-Null
$ws.Range("A1:M92").AutoFilter() | Out-Null

# --- Keep the hidden _FilterDatabase defined name in sync ---------------
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$M`$92"
    }
}

# --- Move the selection to the next empty row, as Excel leaves it after
#     appending rows at the bottom of the sheet -------------------------
$ws.Range("A93").Select() | Out-Null
